# Revert "Add translations to forms"
#
# This undoes a prior edit that had:
#   - added a "table_specific_translations" sheet
#   - replaced short question codes (q65a..q65f) with full English text
#     stored alongside Portuguese/Swahili translations on that sheet
#   - renamed a couple of "display.*" setting keys to "*.text" variants
#   - added English "Male"/"Female"/"Don't Know" display titles for choices
#
# The revert:
#   - restores the full English prompt text directly on the "survey" sheet
#     (instead of the short q65x placeholder codes)
#   - restores the "display.prompt.text" header name on "survey"
#   - restores "display.title.text" header name + Male/Female/Don't Know
#     display titles on "choices"
#   - removes the now-redundant "table_specific_translations" sheet
#   - leaves the active sheet on "choices" (previously "survey")

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet: restore full prompt text -------------------------------
$survey.Range("F1").Value  = "display.prompt.text"
$survey.Range("F3").Value  = "65a. ExtID of the person who died"
$survey.Range("F4").Value  = "65b. Given name of the person who died"
$survey.Range("F5").Value  = "65c. Surname of the person who died"
$survey.Range("F6").Value  = "65d. Gender of the person who died"
$survey.Range("F7").Value  = "65e. Date of death of the person who died"
$survey.Range("F12").Value = "65f. Approximate age of person who died"

# column F now holds much longer strings - widen it to fit
$survey.Columns.Item(6).ColumnWidth = 33.8

# --- choices sheet: restore display.title.text + English titles -----------
$choices.Range("C1").Value = "display.title.text"
$choices.Range("C2").Value = "Male"
$choices.Range("C3").Value = "Female"
$choices.Range("C4").Value = "Don't Know"

# --- remove the table_specific_translations sheet --------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("table_specific_translations").Delete() | Out-Null

# --- restore sheet selections / active sheet --------------------------------
# "survey" was the active tab before; after the revert "choices" is active.
$survey.Activate() | Out-Null
$survey.Range("A3").Select() | Out-Null

$choices.Activate() | Out-Null
$choices.Range("C5").Select() | Out-Null

Write-Output "revert complete"
